# DSM Scheduled Flights vs actual.xlsx
# Append new daily flight rows (2022-10-18 .. 2022-11-09) to the bottom of
# the "Ark1" sheet's data table, continuing the existing
# date / scheduled(B) / actual(C) / on-time-rate(D=C/B) pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, date(serial), scheduled, actual
$data = @(
  @(924, 44852, 58, 54),
  @(925, 44853, 68, 65),
  @(926, 44854, 90, 89),
  @(927, 44855, 71, 71),
  @(928, 44856, 53, 52),
  @(929, 44857, 61, 60),
  @(930, 44858, 71, 70),
  @(931, 44859, 68, 67),
  @(932, 44860, 61, 61),
  @(933, 44861, 86, 82),
  @(934, 44862, 78, 77),
  @(935, 44863, 51, 50),
  @(936, 44864, 55, 54),
  @(937, 44865, 60, 56),
  @(938, 44866, 57, 56),
  @(939, 44867, 73, 68),
  @(940, 44868, 76, 73),
  @(941, 44869, 69, 66),
  @(942, 44870, 49, 46),
  @(943, 44871, 61, 58),
  @(944, 44872, 55, 55),
  @(945, 44873, 71, 70),
  @(946, 44874, 56, 54)
)

# The last existing data row (923) carries the formatting (date / integer /
# percentage styles) we want the new rows to inherit, same as dragging the
# fill handle down in Excel.
$lastRow = 923

foreach ($row in $data) {
  $r = $row[0]

  $src = $ws.Range("A" + $lastRow + ":D" + $lastRow)
  $dst = $ws.Range("A" + $r + ":D" + $r)
  $src.Copy($dst)

  $ws.Range("A" + $r).Value = $row[1]
  $ws.Range("B" + $r).Value = $row[2]
  $ws.Range("C" + $r).Value = $row[3]
  $ws.Range("D" + $r).Formula = "=C" + $r + "/B" + $r

  $lastRow = $r
}

# Match the saved selection/scroll state from the edit: the cursor ends up
# on the newly-filled D column block.
$win = $excel.ActiveWindow
$win.ScrollRow = 918
$win.ScrollColumn = 1
$ws.Range("D921:D946").Select() | Out-Null
